# "Update countries & provincias Spain"
# - Swap the Malaui/Maldivas and Sri Lanka/Suazilandia row labels
#   (their relative order in the source list changed).
# - Refresh the COVID-19 stat counters for the countries whose figures moved.
# - Bump the "last updated" timestamp string in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name corrections (order swapped in the upstream country list) ---
$ws.Cells.Item(107, 1).Value = "Maldivas"
$ws.Cells.Item(108, 1).Value = "Malaui"

$ws.Cells.Item(118, 1).Value = "Suazilandia"
$ws.Cells.Item(119, 1).Value = "Sri Lanka"

# --- Refreshed statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Estados Unidos
$ws.Cells.Item(4, 2).Value = 4836533
$ws.Cells.Item(4, 3).Value = 22886
$ws.Cells.Item(4, 4).Value = 2392119
$ws.Cells.Item(4, 5).Value = 2285817
$ws.Cells.Item(4, 7).Value = 232
$ws.Cells.Item(4, 8).Value = 158597

# India
$ws.Cells.Item(6, 2).Value = 1855331
$ws.Cells.Item(6, 3).Value = 50629
$ws.Cells.Item(6, 4).Value = 1230440
$ws.Cells.Item(6, 5).Value = 585920
$ws.Cells.Item(6, 7).Value = 810
$ws.Cells.Item(6, 8).Value = 38971

# Row 21
$ws.Cells.Item(21, 2).Value = 212158
$ws.Cells.Item(21, 3).Value = 696
$ws.Cells.Item(21, 5).Value = 9330
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 9228

# Row 23
$ws.Cells.Item(23, 2).Value = 191295
$ws.Cells.Item(23, 3).Value = 556
$ws.Cells.Item(23, 5).Value = 79501
$ws.Cells.Item(23, 7).Value = 22
$ws.Cells.Item(23, 8).Value = 30294

# Row 31
$ws.Cells.Item(31, 2).Value = 87041
$ws.Cells.Item(31, 3).Value = 517
$ws.Cells.Item(31, 5).Value = 21930
$ws.Cells.Item(31, 7).Value = 17
$ws.Cells.Item(31, 8).Value = 5767

# Row 52
$ws.Cells.Item(52, 5).Value = 2721
$ws.Cells.Item(52, 7).Value = 2
$ws.Cells.Item(52, 8).Value = 149

# Row 89
$ws.Cells.Item(89, 2).Value = 7948
$ws.Cells.Item(89, 3).Value = 91
$ws.Cells.Item(89, 4).Value = 6767
$ws.Cells.Item(89, 5).Value = 1137
$ws.Cells.Item(89, 7).Value = 1
$ws.Cells.Item(89, 8).Value = 44

# Row 104
$ws.Cells.Item(104, 5).Value = 3154
$ws.Cells.Item(104, 7).Value = 1
$ws.Cells.Item(104, 8).Value = 209

# Row 107 (Maldivas, after relabel)
$ws.Cells.Item(107, 2).Value = 4293
$ws.Cells.Item(107, 3).Value = 129
$ws.Cells.Item(107, 4).Value = 2670
$ws.Cells.Item(107, 5).Value = 1605
$ws.Cells.Item(107, 8).Value = 18

# Row 108 (Malaui, after relabel)
$ws.Cells.Item(108, 2).Value = 4231
$ws.Cells.Item(108, 4).Value = 1919
$ws.Cells.Item(108, 5).Value = 2189
$ws.Cells.Item(108, 8).Value = 123

# Row 112
$ws.Cells.Item(112, 5).Value = 1517
$ws.Cells.Item(112, 7).Value = 4
$ws.Cells.Item(112, 8).Value = 38

# Row 118 (Suazilandia, after relabel)
$ws.Cells.Item(118, 2).Value = 2838
$ws.Cells.Item(118, 3).Value = 63
$ws.Cells.Item(118, 4).Value = 1253
$ws.Cells.Item(118, 5).Value = 1540
$ws.Cells.Item(118, 7).Value = 2
$ws.Cells.Item(118, 8).Value = 45

# Row 119 (Sri Lanka, after relabel)
$ws.Cells.Item(119, 2).Value = 2828
$ws.Cells.Item(119, 3).Value = 5
$ws.Cells.Item(119, 4).Value = 2517
$ws.Cells.Item(119, 5).Value = 300
$ws.Cells.Item(119, 8).Value = 11

# Row 122
$ws.Cells.Item(122, 2).Value = 2543
$ws.Cells.Item(122, 3).Value = 2
$ws.Cells.Item(122, 5).Value = 476

# Row 131
$ws.Cells.Item(131, 4).Value = 676
$ws.Cells.Item(131, 5).Value = 1283

# Row 134
$ws.Cells.Item(134, 2).Value = 1848
$ws.Cells.Item(134, 3).Value = 5
$ws.Cells.Item(134, 5).Value = 406

# Row 138
$ws.Cells.Item(138, 2).Value = 1565
$ws.Cells.Item(138, 3).Value = 4
$ws.Cells.Item(138, 4).Value = 1225

# Row 155
$ws.Cells.Item(155, 2).Value = 847
$ws.Cells.Item(155, 3).Value = 38
$ws.Cells.Item(155, 4).Value = 268
$ws.Cells.Item(155, 5).Value = 533
$ws.Cells.Item(155, 7).Value = 2
$ws.Cells.Item(155, 8).Value = 46

# Row 186
$ws.Cells.Item(186, 2).Value = 121
$ws.Cells.Item(186, 3).Value = 1
$ws.Cells.Item(186, 5).Value = 12

# --- Bump the "last updated" timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Agosto de 2020 a las 21:26"
